# Add the "AddVacancy" worksheet (data-driven test data for the Add Vacancy
# test) right after Sheet1, fill it in with the job-vacancy table, format
# column A as text, size the columns to fit their content and make the new
# sheet the active tab with A5 selected.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "AddVacancy"

# Column A (the numeric "row id") is stored as text, like in the source file.
$ws.Columns.Item(1).NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "job title"
$ws.Range("B1").Value = "vacancy name"
$ws.Range("C1").Value = "hiring manager"
$ws.Range("D1").Value = "number of positions"
$ws.Range("E1").Value = "description"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "test engineer"
$ws.Range("C2").Value = "Kallyani Bhute"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "perform test using selenium in java"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "production"
$ws.Range("C3").Value = "Paul Collings"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "part of team in assembly line"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "driver"
$ws.Range("C4").Value = "Rebecca Harmony"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = "AZ driver"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "developer"
$ws.Range("C5").Value = "Dominic Chase"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "frontend developer"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "analyst"
$ws.Range("C6").Value = "Nathan Elliot"
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = "anylyse software feasibility"

# Size each column to (approximately) fit its longest entry, matching the
# widths Excel's own "best fit" produced on the source workbook.
$ws.Columns.Item(1).ColumnWidth = 6.666666666666667
$ws.Columns.Item(2).ColumnWidth = 11.833333333333332
$ws.Columns.Item(3).ColumnWidth = 15.0
$ws.Columns.Item(4).ColumnWidth = 16.833333333333336
$ws.Columns.Item(5).ColumnWidth = 29.833333333333336

# Make the new sheet the active tab, with A5 selected (mirrors the saved
# selection state in the edited workbook).
$ws.Activate()
$ws.Range("A5").Select()
